# Auto-generated edit script: adds a new weekly price entry (row 7) and
# shifts the existing rows down by one, appending the former last row at the end.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 53: same constant columns as every other data row, with the
#     Volumen / Precios that used to live in row 52 ("bumped" down). ---
$ws.Range("A53").Value = 10
$ws.Range("B53").Value = "Vega Modelo de Temuco"
$ws.Range("C53").Value = "La Araucanía"
$ws.Range("D53").Value = 44425
$ws.Range("E53").Value = 9
$ws.Range("F53").Value = 100112035
$ws.Range("G53").Value = "Bruselas (repollito)"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 30
$ws.Range("K53").Value = 25000
$ws.Range("L53").Value = 25000
$ws.Range("M53").Value = 25000
$ws.Range("N53").Value = "$/malla 10 kilos"
$ws.Range("O53").Value = "Provincia de Quillota"
$ws.Range("P53").Value = 2500
$ws.Range("Q53").Value = 10
$ws.Range("R53").Value = "Hortaliza"
$ws.Range("D53").NumberFormat = $ws.Range("D52").NumberFormat

# --- Shift rows 52 down to 8: row r takes the Fecha/Volumen/Precios that
#     used to belong to row (r-1). Walking bottom-up so no value is lost
#     before it has been copied onward. ---
$ws.Range("D52").Value = 44376
$ws.Range("J52").Value = 45
$ws.Range("K52").Value = 23000
$ws.Range("L52").Value = 23000
$ws.Range("M52").Value = 23000
$ws.Range("P52").Value = 2300
$ws.Range("D51").Value = 44411
$ws.Range("J51").Value = 40
$ws.Range("K51").Value = 25000
$ws.Range("L51").Value = 25000
$ws.Range("M51").Value = 25000
$ws.Range("P51").Value = 2500
$ws.Range("D50").Value = 44410
$ws.Range("J50").Value = 50
$ws.Range("K50").Value = 25000
$ws.Range("L50").Value = 25000
$ws.Range("M50").Value = 25000
$ws.Range("P50").Value = 2500
$ws.Range("D49").Value = 44461
$ws.Range("J49").Value = 40
$ws.Range("K49").Value = 25000
$ws.Range("L49").Value = 25000
$ws.Range("M49").Value = 25000
$ws.Range("P49").Value = 2500
$ws.Range("D48").Value = 44365
$ws.Range("J48").Value = 85
$ws.Range("K48").Value = 22000
$ws.Range("L48").Value = 22000
$ws.Range("M48").Value = 22000
$ws.Range("P48").Value = 2200
$ws.Range("D47").Value = 44356
$ws.Range("J47").Value = 15
$ws.Range("K47").Value = 24000
$ws.Range("L47").Value = 24000
$ws.Range("M47").Value = 24000
$ws.Range("P47").Value = 2400
$ws.Range("D46").Value = 44372
$ws.Range("J46").Value = 20
$ws.Range("K46").Value = 25000
$ws.Range("L46").Value = 25000
$ws.Range("M46").Value = 25000
$ws.Range("P46").Value = 2500
$ws.Range("D45").Value = 44438
$ws.Range("J45").Value = 50
$ws.Range("K45").Value = 25000
$ws.Range("L45").Value = 25000
$ws.Range("M45").Value = 25000
$ws.Range("P45").Value = 2500
$ws.Range("D44").Value = 44397
$ws.Range("J44").Value = 30
$ws.Range("K44").Value = 27000
$ws.Range("L44").Value = 27000
$ws.Range("M44").Value = 27000
$ws.Range("P44").Value = 2700
$ws.Range("D43").Value = 44371
$ws.Range("J43").Value = 50
$ws.Range("K43").Value = 25000
$ws.Range("L43").Value = 25000
$ws.Range("M43").Value = 25000
$ws.Range("P43").Value = 2500
$ws.Range("D42").Value = 44433
$ws.Range("J42").Value = 25
$ws.Range("K42").Value = 25000
$ws.Range("L42").Value = 25000
$ws.Range("M42").Value = 25000
$ws.Range("P42").Value = 2500
$ws.Range("D41").Value = 44386
$ws.Range("J41").Value = 20
$ws.Range("K41").Value = 25000
$ws.Range("L41").Value = 25000
$ws.Range("M41").Value = 25000
$ws.Range("P41").Value = 2500
$ws.Range("D40").Value = 44390
$ws.Range("J40").Value = 15
$ws.Range("K40").Value = 25000
$ws.Range("L40").Value = 25000
$ws.Range("M40").Value = 25000
$ws.Range("P40").Value = 2500
$ws.Range("D39").Value = 44447
$ws.Range("J39").Value = 30
$ws.Range("K39").Value = 27000
$ws.Range("L39").Value = 27000
$ws.Range("M39").Value = 27000
$ws.Range("P39").Value = 2700
$ws.Range("D38").Value = 44389
$ws.Range("J38").Value = 65
$ws.Range("K38").Value = 25000
$ws.Range("L38").Value = 25000
$ws.Range("M38").Value = 25000
$ws.Range("P38").Value = 2500
$ws.Range("D37").Value = 44427
$ws.Range("J37").Value = 40
$ws.Range("K37").Value = 25000
$ws.Range("L37").Value = 25000
$ws.Range("M37").Value = 25000
$ws.Range("P37").Value = 2500
$ws.Range("D36").Value = 44455
$ws.Range("J36").Value = 20
$ws.Range("K36").Value = 25000
$ws.Range("L36").Value = 25000
$ws.Range("M36").Value = 25000
$ws.Range("P36").Value = 2500
$ws.Range("D35").Value = 44350
$ws.Range("J35").Value = 40
$ws.Range("K35").Value = 24000
$ws.Range("L35").Value = 25000
$ws.Range("M35").Value = 24375
$ws.Range("P35").Value = 2438
$ws.Range("D34").Value = 44420
$ws.Range("J34").Value = 55
$ws.Range("K34").Value = 25000
$ws.Range("L34").Value = 25000
$ws.Range("M34").Value = 25000
$ws.Range("P34").Value = 2500
$ws.Range("D33").Value = 44354
$ws.Range("J33").Value = 30
$ws.Range("K33").Value = 24000
$ws.Range("L33").Value = 24000
$ws.Range("M33").Value = 24000
$ws.Range("P33").Value = 2400
$ws.Range("D32").Value = 44448
$ws.Range("J32").Value = 15
$ws.Range("K32").Value = 25000
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = 25000
$ws.Range("P32").Value = 2500
$ws.Range("D31").Value = 44453
$ws.Range("J31").Value = 40
$ws.Range("K31").Value = 27000
$ws.Range("L31").Value = 27000
$ws.Range("M31").Value = 27000
$ws.Range("P31").Value = 2700
$ws.Range("D30").Value = 44392
$ws.Range("J30").Value = 25
$ws.Range("K30").Value = 24000
$ws.Range("L30").Value = 24000
$ws.Range("M30").Value = 24000
$ws.Range("P30").Value = 2400
$ws.Range("D29").Value = 44400
$ws.Range("J29").Value = 12
$ws.Range("K29").Value = 24000
$ws.Range("L29").Value = 24000
$ws.Range("M29").Value = 24000
$ws.Range("P29").Value = 2400
$ws.Range("D28").Value = 44426
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = 25000
$ws.Range("L28").Value = 25000
$ws.Range("M28").Value = 25000
$ws.Range("P28").Value = 2500
$ws.Range("D27").Value = 44405
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 25000
$ws.Range("L27").Value = 25000
$ws.Range("M27").Value = 25000
$ws.Range("P27").Value = 2500
$ws.Range("D26").Value = 44441
$ws.Range("J26").Value = 70
$ws.Range("K26").Value = 25000
$ws.Range("L26").Value = 25000
$ws.Range("M26").Value = 25000
$ws.Range("P26").Value = 2500
$ws.Range("D25").Value = 44452
$ws.Range("J25").Value = 80
$ws.Range("K25").Value = 25000
$ws.Range("L25").Value = 25000
$ws.Range("M25").Value = 25000
$ws.Range("P25").Value = 2500
$ws.Range("D24").Value = 44434
$ws.Range("J24").Value = 55
$ws.Range("K24").Value = 25000
$ws.Range("L24").Value = 25000
$ws.Range("M24").Value = 25000
$ws.Range("P24").Value = 2500
$ws.Range("D23").Value = 44419
$ws.Range("J23").Value = 25
$ws.Range("K23").Value = 25000
$ws.Range("L23").Value = 25000
$ws.Range("M23").Value = 25000
$ws.Range("P23").Value = 2500
$ws.Range("D22").Value = 44417
$ws.Range("J22").Value = 15
$ws.Range("K22").Value = 25000
$ws.Range("L22").Value = 25000
$ws.Range("M22").Value = 25000
$ws.Range("P22").Value = 2500
$ws.Range("D21").Value = 44385
$ws.Range("J21").Value = 80
$ws.Range("K21").Value = 25000
$ws.Range("L21").Value = 25000
$ws.Range("M21").Value = 25000
$ws.Range("P21").Value = 2500
$ws.Range("D20").Value = 44406
$ws.Range("J20").Value = 120
$ws.Range("K20").Value = 24000
$ws.Range("L20").Value = 25000
$ws.Range("M20").Value = 24542
$ws.Range("P20").Value = 2454
$ws.Range("D19").Value = 44446
$ws.Range("J19").Value = 40
$ws.Range("K19").Value = 27000
$ws.Range("L19").Value = 27000
$ws.Range("M19").Value = 27000
$ws.Range("P19").Value = 2700
$ws.Range("D18").Value = 44355
$ws.Range("J18").Value = 25
$ws.Range("K18").Value = 23000
$ws.Range("L18").Value = 24000
$ws.Range("M18").Value = 23400
$ws.Range("P18").Value = 2340
$ws.Range("D17").Value = 44449
$ws.Range("J17").Value = 12
$ws.Range("K17").Value = 25000
$ws.Range("L17").Value = 25000
$ws.Range("M17").Value = 25000
$ws.Range("P17").Value = 2500
$ws.Range("D16").Value = 44435
$ws.Range("J16").Value = 185
$ws.Range("K16").Value = 25000
$ws.Range("L16").Value = 27000
$ws.Range("M16").Value = 25162
$ws.Range("P16").Value = 2516
$ws.Range("D15").Value = 44348
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 25000
$ws.Range("L15").Value = 25000
$ws.Range("M15").Value = 25000
$ws.Range("P15").Value = 2500
$ws.Range("D14").Value = 44431
$ws.Range("J14").Value = 65
$ws.Range("K14").Value = 25000
$ws.Range("L14").Value = 25000
$ws.Range("M14").Value = 25000
$ws.Range("P14").Value = 2500
$ws.Range("D13").Value = 44412
$ws.Range("J13").Value = 50
$ws.Range("K13").Value = 25000
$ws.Range("L13").Value = 25000
$ws.Range("M13").Value = 25000
$ws.Range("P13").Value = 2500
$ws.Range("D12").Value = 44396
$ws.Range("J12").Value = 20
$ws.Range("K12").Value = 25000
$ws.Range("L12").Value = 25000
$ws.Range("M12").Value = 25000
$ws.Range("P12").Value = 2500
$ws.Range("D11").Value = 44349
$ws.Range("J11").Value = 45
$ws.Range("K11").Value = 24000
$ws.Range("L11").Value = 24000
$ws.Range("M11").Value = 24000
$ws.Range("P11").Value = 2400
$ws.Range("D10").Value = 44384
$ws.Range("J10").Value = 40
$ws.Range("K10").Value = 25000
$ws.Range("L10").Value = 25000
$ws.Range("M10").Value = 25000
$ws.Range("P10").Value = 2500
$ws.Range("D9").Value = 44379
$ws.Range("J9").Value = 35
$ws.Range("K9").Value = 22000
$ws.Range("L9").Value = 22000
$ws.Range("M9").Value = 22000
$ws.Range("P9").Value = 2200
$ws.Range("D8").Value = 44421
$ws.Range("J8").Value = 55
$ws.Range("K8").Value = 25000
$ws.Range("L8").Value = 25000
$ws.Range("M8").Value = 25000
$ws.Range("P8").Value = 2500

# --- Row 7: the new weekly entry (Precio mínimo/máximo/ponderado/Kg unchanged) ---
$ws.Range("D7").Value = 44462
$ws.Range("J7").Value = 50

